$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = 'CC'
$ws.Range("C16").Value = '1047457606'
$ws.Range("D16").Value = 'YULIANYS ALVARADO MEDRANO'
$ws.Range("E16").Value = '1803'
$ws.Range("F16").Value = 32000
$ws.Range("G16").Value = 1200000

$ws.Range("B17").Value = 'CC'
$ws.Range("C17").Value = '1047457606'
$ws.Range("D17").Value = 'YULIANYS ALVARADO MEDRANO'
$ws.Range("E17").Value = '1802'
$ws.Range("F17").Value = 48000
$ws.Range("G17").Value = 1200000

$ws.Range("B18").Value = 'CC'
$ws.Range("C18").Value = '1047457606'
$ws.Range("D18").Value = 'YULIANYS ALVARADO MEDRANO'
$ws.Range("E18").Value = '1801'
$ws.Range("F18").Value = 48000
$ws.Range("G18").Value = 1200000

$ws.Range("B19").Value = 'CC'
$ws.Range("C19").Value = '1047457606'
$ws.Range("D19").Value = 'YULIANYS ALVARADO MEDRANO'
$ws.Range("E19").Value = '1712'
$ws.Range("F19").Value = 48000
$ws.Range("G19").Value = 1200000

$ws.Range("B20").Value = 'CC'
$ws.Range("C20").Value = '1047457606'
$ws.Range("D20").Value = 'YULIANYS ALVARADO MEDRANO'
$ws.Range("E20").Value = '1711'
$ws.Range("F20").Value = 48000
$ws.Range("G20").Value = 1200000

$ws.Range("B21").Value = 'CE'
$ws.Range("C21").Value = '700178'
$ws.Range("D21").Value = 'GUSTAVO ADOLFO MOLLEDA BRAVO'
$ws.Range("E21").Value = '1803'
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1500000

$ws.Range("B22").Value = 'CE'
$ws.Range("C22").Value = '700178'
$ws.Range("D22").Value = 'GUSTAVO ADOLFO MOLLEDA BRAVO'
$ws.Range("E22").Value = '1802'
$ws.Range("F22").Value = 60000
$ws.Range("G22").Value = 1500000

$ws.Range("B23").Value = 'CE'
$ws.Range("C23").Value = '700178'
$ws.Range("D23").Value = 'GUSTAVO ADOLFO MOLLEDA BRAVO'
$ws.Range("E23").Value = '1801'
$ws.Range("F23").Value = 60000
$ws.Range("G23").Value = 1500000

$ws.Range("B24").Value = 'CE'
$ws.Range("C24").Value = '700178'
$ws.Range("D24").Value = 'GUSTAVO ADOLFO MOLLEDA BRAVO'
$ws.Range("E24").Value = '1712'
$ws.Range("F24").Value = 60000
$ws.Range("G24").Value = 1500000

$ws.Range("B25").Value = 'CE'
$ws.Range("C25").Value = '700178'
$ws.Range("D25").Value = 'GUSTAVO ADOLFO MOLLEDA BRAVO'
$ws.Range("E25").Value = '1711'
$ws.Range("F25").Value = 60000
$ws.Range("G25").Value = 1500000

$ws.Range("B26").Value = 'CE'
$ws.Range("C26").Value = '562626'
$ws.Range("D26").Value = 'MARIA GABRIELA DE LA GUADALUPE CHIRINOS MORALES'
$ws.Range("E26").Value = '1803'
$ws.Range("F26").Value = 19673
$ws.Range("G26").Value = 737717

$ws.Range("B27").Value = 'CE'
$ws.Range("C27").Value = '562626'
$ws.Range("D27").Value = 'MARIA GABRIELA DE LA GUADALUPE CHIRINOS MORALES'
$ws.Range("E27").Value = '1802'
$ws.Range("F27").Value = 29509
$ws.Range("G27").Value = 737717

$ws.Range("B28").Value = 'CE'
$ws.Range("C28").Value = '562626'
$ws.Range("D28").Value = 'MARIA GABRIELA DE LA GUADALUPE CHIRINOS MORALES'
$ws.Range("E28").Value = '1801'
$ws.Range("F28").Value = 29509
$ws.Range("G28").Value = 737717

$ws.Range("B29").Value = 'CE'
$ws.Range("C29").Value = '562626'
$ws.Range("D29").Value = 'MARIA GABRIELA DE LA GUADALUPE CHIRINOS MORALES'
$ws.Range("E29").Value = '1712'
$ws.Range("F29").Value = 29509
$ws.Range("G29").Value = 737717

$ws.Range("B30").Value = 'CE'
$ws.Range("C30").Value = '562626'
$ws.Range("D30").Value = 'MARIA GABRIELA DE LA GUADALUPE CHIRINOS MORALES'
$ws.Range("E30").Value = '1711'
$ws.Range("F30").Value = 29509
$ws.Range("G30").Value = 737717

$ws.Range("B31").Value = 'CE'
$ws.Range("C31").Value = '562626'
$ws.Range("D31").Value = 'MARIA GABRIELA DE LA GUADALUPE CHIRINOS MORALES'
$ws.Range("E31").Value = '1710'
$ws.Range("F31").Value = 29509
$ws.Range("G31").Value = 737717

$ws.Range("B32").Value = 'CC'
$ws.Range("C32").Value = '79623974'
$ws.Range("D32").Value = 'NELSON GIOVANY ADARME SILVA'
$ws.Range("E32").Value = '1803'
$ws.Range("F32").Value = 29334
$ws.Range("G32").Value = 1100000

$ws.Range("B33").Value = 'CC'
$ws.Range("C33").Value = '79623974'
$ws.Range("D33").Value = 'NELSON GIOVANY ADARME SILVA'
$ws.Range("E33").Value = '1802'
$ws.Range("F33").Value = 44000
$ws.Range("G33").Value = 1100000

$ws.Range("B34").Value = 'CC'
$ws.Range("C34").Value = '79623974'
$ws.Range("D34").Value = 'NELSON GIOVANY ADARME SILVA'
$ws.Range("E34").Value = '1801'
$ws.Range("F34").Value = 44000
$ws.Range("G34").Value = 1100000

$ws.Range("B35").Value = 'CC'
$ws.Range("C35").Value = '79623974'
$ws.Range("D35").Value = 'NELSON GIOVANY ADARME SILVA'
$ws.Range("E35").Value = '1712'
$ws.Range("F35").Value = 44000
$ws.Range("G35").Value = 1100000

$ws.Range("B36").Value = 'CC'
$ws.Range("C36").Value = '79623974'
$ws.Range("D36").Value = 'NELSON GIOVANY ADARME SILVA'
$ws.Range("E36").Value = '1711'
$ws.Range("F36").Value = 44000
$ws.Range("G36").Value = 1100000

$ws.Range("B37").Value = 'CC'
$ws.Range("C37").Value = '1140897176'
$ws.Range("D37").Value = 'PATRICIA BERENICE BONYUET VENENCIA'
$ws.Range("E37").Value = '1803'
$ws.Range("F37").Value = 19673
$ws.Range("G37").Value = 737717

$ws.Range("B38").Value = 'CC'
$ws.Range("C38").Value = '1140897176'
$ws.Range("D38").Value = 'PATRICIA BERENICE BONYUET VENENCIA'
$ws.Range("E38").Value = '1802'
$ws.Range("F38").Value = 29509
$ws.Range("G38").Value = 737717

$ws.Range("B39").Value = 'CC'
$ws.Range("C39").Value = '1140897176'
$ws.Range("D39").Value = 'PATRICIA BERENICE BONYUET VENENCIA'
$ws.Range("E39").Value = '1801'
$ws.Range("F39").Value = 29509
$ws.Range("G39").Value = 737717

$ws.Range("B40").Value = 'CC'
$ws.Range("C40").Value = '1140897176'
$ws.Range("D40").Value = 'PATRICIA BERENICE BONYUET VENENCIA'
$ws.Range("E40").Value = '1712'
$ws.Range("F40").Value = 29509
$ws.Range("G40").Value = 737717

$ws.Range("B41").Value = 'CC'
$ws.Range("C41").Value = '1140897176'
$ws.Range("D41").Value = 'PATRICIA BERENICE BONYUET VENENCIA'
$ws.Range("E41").Value = '1711'
$ws.Range("F41").Value = 29509
$ws.Range("G41").Value = 737717
